$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Orthodontics", "Unlikely", "Service recieved was adaquete but staff seemed like they dont care at all"),
    @("Theatre Treatment Suite Implants", "Extremely Unlikely", "Service recieved was adaquete but staff seemed like they dont care at all"),
    @("Gynaecology", "Unlikely", "Doctors are patronising and made me feel bad"),
    @("Dermatology", "Extremely Unlikely", "Waited for long time for poor service"),
    @("Special Care Baby Unit", "Unlikely", "Felt as if i was not a priority"),
    @("Rehab Services", "Unlikely", "Waited too long to find a parking spot"),
    @("Day Surgery", "Extremely Unlikely", "Felt as if i was not a priority"),
    @("Radiology", "Unlikely", "doctors dont seem to care about me, felt ignored"),
    @("A&E", "Extremely Unlikely", "doctors dont seem to care about me, felt ignored"),
    @("Sitwell", "Unlikely", "Long wait times"),
    @("Sitwell", "Unlikely", "I waited for a long time to be turned away"),
    @("A&E", "Unlikely", "Clenliness isn't the best but otherwise okay"),
    @("Bone Health", "Extremely Unlikely", "Food was terrible"),
    @("Labour and Delivery Suite", "Unlikely", "staff tried to deal with me quickly rather than correctly. Not appropriate and i shouldve have been taken care of better. Would not recommend."),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value2 = $data[$i][0]
    $ws.Cells.Item($row, 2).Value2 = $data[$i][1]
    $ws.Cells.Item($row, 3).Value2 = $data[$i][2]
}

$ws.Rows("16:16").Select()
